$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.678.80"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.530.67"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.33"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.483"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "1.750.29"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "1.521.78"
$ws.Range("E13").Value = "  -2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.66"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.504"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "26.674.36"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.06"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "211.87"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "0.0₃0680"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.18"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.99"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.04"
$ws.Range("E23").Value = "  -3.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.68"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.54"
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.78"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.09"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0451"
$ws.Range("E31").Value = "  -2.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").Value = "1.357.33"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.933"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0162"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.520"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.68"
$ws.Range("E42").Value = "  +5.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.991"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.73"
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.27"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").Value = "1.663.77"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.29"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0503"
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").Value = "0.0₇0971"
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0940"
$ws.Range("E51").Value = "  -0.83%  "
